$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B.. to the right).
# Excel copies the formatting of the column to the left (A) into the new column.
$ws.Range("B1").EntireColumn.Insert()

# Populate the new column's header and sample row with the second "codigo" fields.
$ws.Range("B1").Value = "Código2"
$ws.Range("B2").Value = "cod02"

# Update the hidden _FilterDatabase defined name so it spans the new column layout.
$names = $ws.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "plantilla!_FilterDatabase") {
        $n.RefersTo = "=plantilla!`$C`$1:`$Q`$1"
    }
}

# Restore the active selection to match the authored state.
$ws.Range("B3").Select()
